$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.042.75"
$ws.Range("E2").Value = "  -0.47%  "
$ws.Range("D3").Value = "1.781.63"
$ws.Range("E3").Value = "  -2.54%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "225.49"
$ws.Range("E5").Value = "  +0.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.549"
$ws.Range("E6").Value = "  -1.66%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.41"
$ws.Range("E8").Value = "  +1.34%  "
$ws.Range("E9").Value = "  -2.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0709"
$ws.Range("E10").Value = "  -2.41%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0935"
$ws.Range("E11").Value = "  +0.61%  "
$ws.Range("D12").Value = "2.040.06"
$ws.Range("E12").Value = "  -2.48%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.95"
$ws.Range("E13").Value = "  +1.02%  "
$ws.Range("D14").Value = "1.786.57"
$ws.Range("E14").Value = "  -2.21%  "
$ws.Range("D15").Value = "33.992.18"
$ws.Range("E15").Value = "  -0.76%  "
$ws.Range("E16").Value = "  -4.19%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.13"
$ws.Range("E17").Value = "  -5.05%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.71"
$ws.Range("E18").Value = "  -3.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "243.81"
$ws.Range("E19").Value = "  -3.06%  "
$ws.Range("D20").Value = "0.0₃0784"
$ws.Range("E20").Value = "  -1.38%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.00"
$ws.Range("E21").Value = "  +0.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.69"
$ws.Range("E22").Value = "  -4.48%  "
$ws.Range("E23").Value = "  -4.83%  "
$ws.Range("E24").Value = "  -4.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "160.02"
$ws.Range("E25").Value = "  -0.40%  "
$ws.Range("E26").Value = "  -2.56%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.02"
$ws.Range("E27").Value = "  -3.60%  "
$ws.Range("E28").Value = "  -2.62%  "
$ws.Range("E30").Value = "  +0.87%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0511"
$ws.Range("E31").Value = "  -4.34%  "
$ws.Range("E32").Value = "  -4.32%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.50"
$ws.Range("E34").Value = "  -5.25%  "
$ws.Range("D35").Value = "1.389.74"
$ws.Range("E35").Value = "  -3.51%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.645"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.05"
$ws.Range("E37").Value = "  -2.11%  "
$ws.Range("E38").Value = "  -2.02%  "
$ws.Range("E39").Value = "  +0.18%  "
$ws.Range("E40").Value = "  +1.13%  "
$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.69"
$ws.Range("E41").Value = "  -2.87%  "
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.910"
$ws.Range("E42").Value = "  -5.84%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "77.83"
$ws.Range("E43").Value = "  -4.94%  "
$ws.Range("D44").Value = "0.0₆0141"
$ws.Range("E44").Value = "  +14.25%  "
$ws.Range("E45").Value = "  +2.14%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.46"
$ws.Range("E46").Value = "  +4.61%  "
$ws.Range("B47").Value = "Kaspa"
$ws.Range("C47").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0496"
$ws.Range("E47").Value = "  -0.46%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "107.57"
$ws.Range("E48").Value = "  -0.09%  "
$ws.Range("E49").Value = "  -4.59%  "
$ws.Range("D50").Value = "1.940.19"
$ws.Range("E50").Value = "  -2.46%  "
$ws.Range("E51").Value = "  +0.16%  "
